$d = $word.ActiveDocument

# 1. Remove the _GoBack bookmark from its original location (first paragraph).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2. Split the paragraph that currently reads:
#    "...technologies.  To be used in Canada, ... working our way through:"
#    into two separate paragraphs, breaking right after the two trailing
#    spaces following "technologies."
$r = $d.Content
$found = $r.Find.Execute("technologies.  To be used in Canada", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $breakRange = $d.Range($r.Start, $r.Start + 15)
    $breakRange.Collapse(0)
    $breakRange.InsertParagraphAfter()
}

# 3. Re-add the _GoBack bookmark, now spanning the final paragraph of the
#    document (the "FHIR has been slower to penetrate..." paragraph).
$last = $d.Paragraphs.Last
$d.Bookmarks.Add("_GoBack", $last.Range)
